$d = $word.ActiveDocument

function Find-ParagraphIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Insert-ItalicParagraphAfter($searchText, $newText) {
    $idx = Find-ParagraphIndex $searchText
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $newIdx = $idx + 1
    $p2 = $d.Paragraphs($newIdx)
    $r = $p2.Range
    $r.Text = $newText
    $r2 = $d.Range($r.Start, $r.Start + $newText.Length)
    $r2.Font.Italic = 1
}

# 1) Update activation date
$null = $d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ativação: 01/01/2023", 2)

# 2) Add English translation after "Objetivos" paragraph text
Insert-ItalicParagraphAfter `
    "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte." `
    "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."

# 3) Add English translation after "Programa resumido" paragraph text
Insert-ItalicParagraphAfter `
    "A definir, de acordo com o tópico programado." `
    "To be defined, according to the programmed topic."

# 4) Add English translation after "Programa" paragraph text
Insert-ItalicParagraphAfter `
    "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação." `
    "The content of this elective course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
